$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Data": add two new leading rows (2023, 2022) above the existing
# table, and twenty new trailing rows (2004 down to 1985) below it, and
# refresh the "Valor" figures for the whole series.
# ---------------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Data")

# Make room for the two new rows (2023 and 2022) right after the header;
# the previously-existing rows 2-18 (2021..2005) shift down to 4-20.
$ws.Rows("2:3").Insert()

$ws.Range("A2:A3").NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "2023"
$ws.Cells.Item(2, 2).Value = 24.9
$ws.Cells.Item(3, 1).Value = "2022"
$ws.Cells.Item(3, 2).Value = 24.5

# Append twenty new rows (2004 down to 1985) after the old last row, which
# is now row 20.
$years = @("2004","2003","2002","2001","2000","1999","1998","1997","1996","1995","1994","1993","1992","1991","1990","1989","1988","1987","1986","1985")
$values = @(16.3,15.7,15.5,16,16.2,16.2,16,15.8,16.1,16.7,17,15.4,15.7,17.3,17.3,19.1,17.2,16.5,17.7,16.5)

$ws.Range("A21:A40").NumberFormat = "@"
for ($i = 0; $i -lt $years.Length; $i++) {
    $r = 21 + $i
    $ws.Cells.Item($r, 1).Value = $years[$i]
    $ws.Cells.Item($r, 2).Value = $values[$i]
}

# ---------------------------------------------------------------------------
# Sheet "Metadata": blank-out A1 to a single space, and insert a new
# "actualizacion" / "Julio 2025" row before the "cita" row.
# ---------------------------------------------------------------------------
$wsMeta = $wb.Worksheets.Item("Metadata")

$wsMeta.Cells.Item(1, 1).Value = " "

$wsMeta.Rows(9).Insert()
$wsMeta.Cells.Item(9, 1).Value = "actualizacion"
$wsMeta.Cells.Item(9, 2).Value = "Julio 2025"
